$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some numeric-looking price values must be forced to remain as plain text
# (matching the workbook author data, which stores all values as strings),
# otherwise Excel auto-converts them into floating point numbers.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "25.894.72"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "1.638.76"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("D5").Value = "214.75"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").Value = "0.5041"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "0.2573"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").Value = "0.06393"
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("D10").Value = "19.63"
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").Value = "0.07790"
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "4.285"
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.649.65"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").Value = "0.5435"
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("D15").Value = "0.0₅7874"
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("D16").Value = "65.08"
$ws.Range("E16").Value = "  +2.50%  "
$ws.Range("D17").Value = "25.951.73"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "1.006"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").Value = "197.84"
$ws.Range("E19").Value = "  -2.58%  "
$ws.Range("D20").Value = "4.388"
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("D21").Value = "9.958"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "6.020"
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("D23").Value = "1.005"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "1.868"
$ws.Range("E24").Value = "  -3.54%  "
$ws.Range("D25").Value = "140.23"
$ws.Range("E25").Value = "  -0.96%  "
$ws.Range("D26").Value = "0.1140"
$ws.Range("D27").Value = "6.863"
$ws.Range("E27").Value = "  +2.34%  "
$ws.Range("D28").Value = "15.72"
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("D29").Value = "1.239"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "0.05004"
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("D31").Value = "3.260"
$ws.Range("D32").Value = "3.192"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("D33").Value = "1.531"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").Value = "2.360"
$ws.Range("E34").Value = "  +0.45%  "
$ws.Range("D35").Value = "0.8942"
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").Value = "2.595"
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("D37").Value = "1.138.61"
$ws.Range("E37").Value = "  -2.87%  "
$ws.Range("D38").Value = "0.5523"
$ws.Range("E38").Value = "  -0.71%  "
$ws.Range("D39").Value = "0.01554"
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("D40").Value = "1.004"
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").Value = "5.701"
$ws.Range("E41").Value = "  +1.03%  "
$ws.Range("D42").Value = "0.8145"
$ws.Range("E42").Value = "  +1.24%  "
$ws.Range("D43").Value = "99.34"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("E44").Value = "  +11.34%  "
$ws.Range("D45").Value = "1.777.29"
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").Value = "0.4536"
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("D47").Value = "55.15"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("D48").Value = "1.004"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").Value = "0.05079"
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.09546"
$ws.Range("E50").Value = "  +3.24%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "1.004"
$ws.Range("E51").Value = "  +0.02%  "
